# "add header and background"
# The "Our Service" shared string (used by Header/C4) is renamed to
# "Our Services" (the old string becomes unused and is dropped; the new
# text is appended to the shared-string table), and the active selection
# on the Header sheet moves from C3 to F10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Header")

$ws.Range("C4").Value = "Our Services"

[void]$ws.Range("F10").Select()
